# Applies the invoice test-data update described by the commit:
#  - Memo_Verification_details sheet (row 2 / the single data row):
#      Invoice_number   : TESTINV71161 -> TESTINV97479
#      PO Number        : 4500000904   -> 4500000891
#      PO Quantity      : 4.00         -> 1.00
#      Revised Tax Code : V0           -> KG
#  - Memo_invoice_Details sheet (row 2 / the single data row):
#      Invoice_number       : TESTINV71161 -> TESTINV97479
#      Invoice_Date         : 2024-03-28   -> 2024-03-27
#      Base_Amount          : 4.0          -> 1
#      IGST                 : 0            -> 0.18
#      Total_Invoice_Amount : 4.0          -> 1.18

$wb = $excel.ActiveWorkbook

$verification = $wb.Worksheets.Item("Memo_Verification_details")

$verification.Range("B2").NumberFormat = "@"
$verification.Range("B2").Value = "TESTINV97479"

$verification.Range("C2").NumberFormat = "@"
$verification.Range("C2").Value = "4500000891"

$verification.Range("E2").NumberFormat = "@"
$verification.Range("E2").Value = "1.00"

$verification.Range("K2").NumberFormat = "@"
$verification.Range("K2").Value = "KG"

$invoiceDetails = $wb.Worksheets.Item("Memo_invoice_Details")

$invoiceDetails.Range("B2").NumberFormat = "@"
$invoiceDetails.Range("B2").Value = "TESTINV97479"

$invoiceDetails.Range("C2").NumberFormat = "@"
$invoiceDetails.Range("C2").Value = "2024-03-27"

$invoiceDetails.Range("E2").NumberFormat = "@"
$invoiceDetails.Range("E2").Value = "1"

$invoiceDetails.Range("I2").NumberFormat = "@"
$invoiceDetails.Range("I2").Value = "0.18"

$invoiceDetails.Range("O2").NumberFormat = "@"
$invoiceDetails.Range("O2").Value = "1.18"
